$d = $word.ActiveDocument

# 1) "Thiện Nhượng" -> "Tuấn Khanh"  (signer name change, table cell)
$d.Content.Find.Execute("Thiện Nhượng", $true, $false, $false, $false, $false, $true, 1, $false, "Tuấn Khanh", 2)

# 2) Remove ", PTGĐ. Vũ Tuấn Khanh" after "Nguyễn Văn Nam"
$d.Content.Find.Execute(", PTGĐ. Vũ Tuấn Khanh", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3) "CUVT-HCM" -> full name
$d.Content.Find.Execute("CUVT-HCM", $true, $false, $false, $false, $false, $true, 1, $false, "Trung tâm cung ứng vật tư - Viễn thông thành phố Hồ Chí Minh", 2)

# 4) "THIỆN NHƯỢNG" -> "TUẤN KHANH" (bold signature block)
$d.Content.Find.Execute("THIỆN NHƯỢNG", $true, $false, $false, $false, $false, $true, 1, $false, "TUẤN KHANH", 2)
